$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.191.52"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "1.939.78"
$ws.Range("E3").Value = "  -4.50%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'242.04"
$ws.Range("E5").Value = "  -3.01%  "
$ws.Range("D6").Value = "'0.606"
$ws.Range("E6").Value = "  -5.18%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'56.10"
$ws.Range("E8").Value = "  -11.04%  "
$ws.Range("D9").Value = "'0.363"
$ws.Range("E9").Value = "  -7.70%  "
$ws.Range("D10").Value = "'55.31"
$ws.Range("E10").Value = "  -4.64%  "
$ws.Range("D11").Value = "'0.0826"
$ws.Range("E11").Value = "  +4.31%  "
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").Value = "'0.824"
$ws.Range("E13").Value = "  -7.63%  "
$ws.Range("D14").Value = "2.221.81"
$ws.Range("E14").Value = "  -4.64%  "
$ws.Range("D15").Value = "'20.99"
$ws.Range("E15").Value = "  -10.53%  "
$ws.Range("D16").Value = "'13.29"
$ws.Range("E16").Value = "  -7.57%  "
$ws.Range("D17").Value = "'5.20"
$ws.Range("D18").Value = "1.932.52"
$ws.Range("E18").Value = "  -4.65%  "
$ws.Range("D19").Value = "36.040.49"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("D20").Value = "'69.63"
$ws.Range("E20").Value = "  -4.06%  "
$ws.Range("D21").Value = "0.0₃0864"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").Value = "'228.01"
$ws.Range("E22").Value = "  -3.82%  "
$ws.Range("D23").Value = "'4.97"
$ws.Range("E23").Value = "  -7.97%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("D26").Value = "'2.26"
$ws.Range("E26").Value = "  -3.74%  "
$ws.Range("D27").Value = "'9.33"
$ws.Range("E27").Value = "  -6.13%  "
$ws.Range("D28").Value = "'162.87"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").Value = "'19.22"
$ws.Range("E29").Value = "  -5.65%  "
$ws.Range("D30").Value = "'0.116"
$ws.Range("E30").Value = "  -16.91%  "
$ws.Range("D31").Value = "'0.117"
$ws.Range("D32").Value = "'1.14"
$ws.Range("E32").Value = "  -4.94%  "
$ws.Range("D33").Value = "'4.68"
$ws.Range("E33").Value = "  -7.84%  "
$ws.Range("D34").Value = "'0.0623"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "'4.27"
$ws.Range("E35").Value = "  -5.90%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").Value = "'6.02"
$ws.Range("E37").Value = "  -8.31%  "
$ws.Range("D38").Value = "'1.79"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").Value = "'2.14"
$ws.Range("E39").Value = "  -10.92%  "
$ws.Range("D40").Value = "'2.84"
$ws.Range("E40").Value = "  -11.94%  "
$ws.Range("D41").Value = "'0.0967"
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("D43").Value = "'1.17"
$ws.Range("E43").Value = "  -7.54%  "
$ws.Range("D44").Value = "'0.0208"
$ws.Range("E44").Value = "  -3.96%  "
$ws.Range("D45").Value = "'15.56"
$ws.Range("E45").Value = "  -8.72%  "

# Rows 46-51 were re-sorted; re-assign Coin/Link/Price/Volume for the new order
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.343.54"
$ws.Range("E46").Value = "  -1.66%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'1.03"
$ws.Range("E47").Value = "  -9.90%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.32"
$ws.Range("E48").Value = "  -5.38%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'87.40"
$ws.Range("E49").Value = "  -7.10%  "

$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "'46.43"
$ws.Range("E50").Value = "  +2.23%  "

$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").Value = "'2.81"
$ws.Range("E51").Value = "  -3.23%  "

